# Salesforce Training schedule.xlsx -- Power Query "Group By" workflow
#
# Original author flow (reconstructed from the OOXML diff):
#  1. Selected the Topic + Duration (Hours) columns on Sheet3 and loaded
#     them into Power Query, which produced a connection-only query and,
#     when closed & loaded as a table, a new sheet "Sheet4" holding a
#     straight copy of those two columns as a table named "Table1".
#  2. Grouped that query by Topic with a Count aggregation and loaded the
#     result as a new table ("Table1_1") on a new sheet named "Table1",
#     with a hidden ExternalData_1 defined name scoped to that sheet.
#  3. Sheet3 itself ends up no longer the active tab; "Table1" becomes the
#     active / selected tab instead.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# 1. Read the source data (Topic = column E, Duration (Hours) = column G)
#    off Sheet3, rows 2.. until the Topic column runs dry.
# ---------------------------------------------------------------------
$topics = @()
$durations = @()
$r = 2
while ($ws3.Cells.Item($r, 5).Value2 -ne $null -and $ws3.Cells.Item($r, 5).Value2 -ne "") {
    $topics += $ws3.Cells.Item($r, 5).Value2
    $durations += $ws3.Cells.Item($r, 7).Value2
    $r = $r + 1
}

# Group by Topic (first-seen order), counting occurrences -> mirrors the
# Power Query "Group By" step (Count of rows per Topic).
$groupNames = @()
$groupCounts = @{}
for ($i = 0; $i -lt $topics.Length; $i++) {
    $t = $topics[$i]
    if ($groupCounts.ContainsKey($t)) {
        $groupCounts[$t] = $groupCounts[$t] + 1
    } else {
        $groupCounts[$t] = 1
        $groupNames += $t
    }
}

# ---------------------------------------------------------------------
# 2. Create the two new sheets. "Sheet4" is minted first internally (it
#    carries the lower sheetId, 4), then "Table1" is inserted right after
#    Sheet3 -- ahead of Sheet4 -- giving the final tab order
#    Sheet3, Table1, Sheet4 while Table1 keeps the higher sheetId (5).
#
#    NOTE: worksheet references returned by Worksheets.Add() are
#    positional anchors that can get rebound by a later Add() at/around
#    the same slot, so every sheet is re-fetched by name (via
#    Worksheets.Item) right before it is actually written to, below.
# ---------------------------------------------------------------------
$tmp4 = $wb.Worksheets.Add($null, $ws3)
$tmp4.Name = "Sheet4"

$tmpT1 = $wb.Worksheets.Add($null, $ws3)
$tmpT1.Name = "Table1"

$wsSheet4 = $wb.Worksheets.Item("Sheet4")
$wsTable1 = $wb.Worksheets.Item("Table1")

# ---------------------------------------------------------------------
# 3. "Sheet4" -- plain copy of Topic / Duration (Hours) columns, loaded
#    as table "Table1" (query-sourced, but a normal table on export).
# ---------------------------------------------------------------------
$wsSheet4.Cells.Item(1, 1).Value = "Topic"
$wsSheet4.Cells.Item(1, 2).Value = " Duration (Hours)"
for ($i = 0; $i -lt $topics.Length; $i++) {
    $wsSheet4.Cells.Item($i + 2, 1).Value = $topics[$i]
    $wsSheet4.Cells.Item($i + 2, 2).Value = $durations[$i]
}

$loSheet4 = $wsSheet4.ListObjects.Add(1, $wsSheet4.Range("A1:B" + ($topics.Length + 1)), $null, 1)
$loSheet4.Name = "Table1"
$loSheet4.TableStyle = "TableStyleMedium2"

$wsSheet4.Columns.Item(1).ColumnWidth = 27.109375
$wsSheet4.Columns.Item(2).ColumnWidth = 17.33203125

# ---------------------------------------------------------------------
# 4. "Table1" -- grouped Topic / Count result table "Table1_1".
# ---------------------------------------------------------------------
$wsTable1.Cells.Item(1, 1).Value = "Topic"
$wsTable1.Cells.Item(1, 2).Value = "Count"
for ($i = 0; $i -lt $groupNames.Length; $i++) {
    $wsTable1.Cells.Item($i + 2, 1).Value = $groupNames[$i]
    $wsTable1.Cells.Item($i + 2, 2).Value = $groupCounts[$groupNames[$i]]
}

$loTable1 = $wsTable1.ListObjects.Add(1, $wsTable1.Range("A1:B" + ($groupNames.Length + 1)), $null, 1)
$loTable1.Name = "Table1_1"
$loTable1.TableStyle = "TableStyleMedium7"

$wsTable1.Columns.Item(1).ColumnWidth = 27.109375
$wsTable1.Columns.Item(2).ColumnWidth = 8.33203125

# Hidden per-sheet defined name Power Query stamps next to a loaded table.
$extName = $wsTable1.Names.Add("ExternalData_1", "=Table1!`$A`$1:`$B`$" + ($groupNames.Length + 1))
$extName.Visible = $false

# ---------------------------------------------------------------------
# 5. View state: Sheet3 is no longer the selected tab; Table1 is.
# ---------------------------------------------------------------------
$ws3.Range("G1").Select()
$wsTable1.Range("G19").Select()
$wsTable1.Activate()

Write-Host "done"
